$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("TotalObjectTranslation - T1")
$ws.Range("B7").Value = 1.07
$ws.Range("B12").Value = 1.09
$ws.Range("B13").Value = 1.65
$ws.Range("B14").Value = 1.38
$ws.Range("B15").Value = 1.14
$ws.Range("B18").Value = 1.72
$ws.Range("B22").Value = 0.9399999999999999

$ws = $wb.Worksheets.Item("TotalObjectTranslation - T2")
$ws.Range("B2").Value = 1.76
$ws.Range("B3").Value = 1.45
$ws.Range("B4").Value = 2.74
$ws.Range("B5").Value = 10.42
$ws.Range("B6").Value = 1.28
$ws.Range("B7").Value = 1.72
$ws.Range("B8").Value = 2.86
$ws.Range("B9").Value = 1.87
$ws.Range("B10").Value = 1.46
$ws.Range("B11").Value = 1.57
$ws.Range("B12").Value = 2.59
$ws.Range("B13").Value = 1.86
$ws.Range("B14").Value = 2.45
$ws.Range("B15").Value = 2.04
$ws.Range("B16").Value = 1.49
$ws.Range("B17").Value = 1.7
$ws.Range("B18").Value = 1.7
$ws.Range("B19").Value = 1.85
$ws.Range("B20").Value = 1.34
$ws.Range("B21").Value = 1.53
$ws.Range("B22").Value = 1.4
$ws.Range("B23").Value = 1.6
$ws.Range("B24").Value = 1.46
$ws.Range("B25").Value = 1.78
$ws.Range("B26").Value = 1.47
$ws.Range("B27").Value = 1.45

$ws = $wb.Worksheets.Item("TotalObjectTranslation - T3")
$ws.Range("B2").Value = 5.05
$ws.Range("B3").Value = 3.87
$ws.Range("B4").Value = 5.22
$ws.Range("B5").Value = 4.32
$ws.Range("B6").Value = 4.91
$ws.Range("B7").Value = 5.08
$ws.Range("B8").Value = 3.67
$ws.Range("B9").Value = 3.72
$ws.Range("B10").Value = 4.37
$ws.Range("B11").Value = 5.88
$ws.Range("B12").Value = 9.19
$ws.Range("B13").Value = 4.54
$ws.Range("B14").Value = 5.62
$ws.Range("B15").Value = 4.08
$ws.Range("B16").Value = 6.72
$ws.Range("B17").Value = 4.46
$ws.Range("B18").Value = 3.92
$ws.Range("B19").Value = 5.220000000000001
$ws.Range("B20").Value = 3.58
$ws.Range("B21").Value = 4.15
$ws.Range("B22").Value = 5.32
$ws.Range("B23").Value = 5.41
$ws.Range("B24").Value = 4.48
$ws.Range("B25").Value = 4.140000000000001
$ws.Range("B26").Value = 5.890000000000001
$ws.Range("B27").Value = 4.08

$ws = $wb.Worksheets.Item("TotalObjectTranslation - T4")
$ws.Range("B2").Value = 3.63
$ws.Range("B3").Value = 3.46
$ws.Range("B4").Value = 3.77
$ws.Range("B5").Value = 3.76
$ws.Range("B6").Value = 3.56
$ws.Range("B7").Value = 3.55
$ws.Range("B8").Value = 4.02
$ws.Range("B9").Value = 3.46
$ws.Range("B10").Value = 4.539999999999999
$ws.Range("B11").Value = 3.84
$ws.Range("B12").Value = 3.95
$ws.Range("B13").Value = 3.48
$ws.Range("B14").Value = 4.79
$ws.Range("B15").Value = 3.55
$ws.Range("B16").Value = 3.68
$ws.Range("B17").Value = 3.8
$ws.Range("B18").Value = 5.12
$ws.Range("B19").Value = 3.58
$ws.Range("B20").Value = 3.49
$ws.Range("B21").Value = 3.62
$ws.Range("B22").Value = 3.45
$ws.Range("B23").Value = 5.100000000000001
$ws.Range("B24").Value = 3.47
$ws.Range("B25").Value = 4.03
$ws.Range("B26").Value = 3.64
$ws.Range("B27").Value = 3.63

$ws = $wb.Worksheets.Item("TotalObjectTranslation - T5")
$ws.Range("B2").Value = 70.32000000000001
$ws.Range("B3").Value = 64.92
$ws.Range("B4").Value = 46.59
$ws.Range("B5").Value = 139.53
$ws.Range("B6").Value = 67.65000000000001
$ws.Range("B7").Value = 48.51000000000001
$ws.Range("B8").Value = 62.39
$ws.Range("B9").Value = 35.65
$ws.Range("B10").Value = 223.53
$ws.Range("B11").Value = 59.95
$ws.Range("B12").Value = 40
$ws.Range("B13").Value = 83.93000000000001
$ws.Range("B14").Value = 62.72
$ws.Range("B15").Value = 40.26
$ws.Range("B16").Value = 37.02
$ws.Range("B17").Value = 97.48
$ws.Range("B18").Value = 47.7
$ws.Range("B19").Value = 38.56
$ws.Range("B20").Value = 85.34
$ws.Range("B21").Value = 50.53
$ws.Range("B22").Value = 148.52
$ws.Range("B23").Value = 48.5
$ws.Range("B24").Value = 42.75
$ws.Range("B25").Value = 125.16
$ws.Range("B26").Value = 51.91
$ws.Range("B27").Value = 51.18

$ws = $wb.Worksheets.Item("TotalObjectTranslation - T6")
$ws.Range("B2").Value = 97.33000000000001
$ws.Range("B3").Value = 63.09
$ws.Range("B4").Value = 68.95
$ws.Range("B5").Value = 76.08
$ws.Range("B6").Value = 71.59999999999999
$ws.Range("B7").Value = 72.72
$ws.Range("B8").Value = 89.23
$ws.Range("B9").Value = 88.93000000000001
$ws.Range("B10").Value = 105.71
$ws.Range("B11").Value = 74.25
$ws.Range("B12").Value = 67.90000000000001
$ws.Range("B13").Value = 89.22
$ws.Range("B14").Value = 112.75
$ws.Range("B15").Value = 89.84
$ws.Range("B16").Value = 81.74000000000001
$ws.Range("B17").Value = 70.83
$ws.Range("B18").Value = 128.95
$ws.Range("B19").Value = 109.59
$ws.Range("B20").Value = 74
$ws.Range("B21").Value = 75.42999999999999
$ws.Range("B22").Value = 86.38
$ws.Range("B23").Value = 79.34
$ws.Range("B24").Value = 69.05
$ws.Range("B25").Value = 62.68
$ws.Range("B26").Value = 69.43000000000001
$ws.Range("B27").Value = 67.75

$ws = $wb.Worksheets.Item("TotalObjectTranslation - T7")
$ws.Range("B4").Value = 49.64
$ws.Range("B12").Value = 52.37
$ws.Range("B14").Value = 62.34
$ws.Range("B15").Value = 189.22
$ws.Range("B16").Value = 62.63
$ws.Range("B20").Value = 68.72999999999999
$ws.Range("B23").Value = 102.18
$ws.Range("B26").Value = 54.5

$ws = $wb.Worksheets.Item("TotalObjectTranslation - T8")
$ws.Range("B2").Value = 139.11
$ws.Range("B3").Value = 143.54
$ws.Range("B4").Value = 149.57
$ws.Range("B5").Value = 164.23
$ws.Range("B6").Value = 124.89
$ws.Range("B7").Value = 121.02
$ws.Range("B8").Value = 130.66
$ws.Range("B9").Value = 116.89
$ws.Range("B10").Value = 268.47
$ws.Range("B11").Value = 277.68
$ws.Range("B12").Value = 103.74
$ws.Range("B13").Value = 140.07
$ws.Range("B14").Value = 126
$ws.Range("B15").Value = 154.59
$ws.Range("B16").Value = 187.21
$ws.Range("B17").Value = 151.26
$ws.Range("B18").Value = 169.35
$ws.Range("B19").Value = 114.17
$ws.Range("B20").Value = 105.87
$ws.Range("B21").Value = 120.51
$ws.Range("B22").Value = 103.64
$ws.Range("B23").Value = 138.04
$ws.Range("B24").Value = 97.45
$ws.Range("B25").Value = 94.36999999999999
$ws.Range("B26").Value = 146.26
$ws.Range("B27").Value = 233.03
